$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.428.82'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '2.252.01'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.45'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.632'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '76.15'
$ws.Range('E7').Value = '  -0.86%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.621'
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '44.28'
$ws.Range('E10').Value = '  +7.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0951'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.27'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.102'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.63'
$ws.Range('E14').Value = '  -1.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.861'
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('D16').Value = '2.250.86'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').Value = '42.328.24'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('E18').Value = '  +4.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.20'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.27'
$ws.Range('E20').Value = '  +1.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.27'
$ws.Range('E21').Value = '  +3.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '232.06'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.90'
$ws.Range('E23').Value = '  +22.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.50'
$ws.Range('E25').Value = '  +2.68%  '
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.49'
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.71'
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0824'
$ws.Range('E31').Value = '  -3.47%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.120'
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '30.86'
$ws.Range('E33').Value = '  -6.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.35'
$ws.Range('E34').Value = '  +10.22%  '
$ws.Range('E35').Value = '  -0.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.55'
$ws.Range('E36').Value = '  -1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0316'
$ws.Range('E37').Value = '  +6.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '14.06'
$ws.Range('E38').Value = '  +5.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.19'
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.81'
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '63.99'
$ws.Range('E41').Value = '  +6.27%  '
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '107.95'
$ws.Range('E43').Value = '  -5.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.80'
$ws.Range('E44').Value = '  -0.40%  '
$ws.Range('E45').Value = '  +2.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.997'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.36'
$ws.Range('E49').Value = '  +3.36%  '
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.14'
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.71'
$ws.Range('E51').Value = '  +0.98%  '
